# This script rotates the data of rows 5, 6 and 7 (the three data records
# following the header row) on the active worksheet:
#   new row 5 <- old row 7
#   new row 6 <- old row 5
#   new row 7 <- old row 6
#
# Every cell in the A:AY range for those rows is copied so that cells
# which are populated in one row but empty in another end up in the right
# place after the rotation.
#
# Columns I, Y and AA store plain text that merely *looks* numeric/date-
# like ("1", "2", "5", "2023-08-31"). Excel's automatic type detection
# would silently turn such a value into a real number/date when it is
# assigned with .Value, so for those columns the destination cell is
# temporarily switched to Text number format before the assignment (and
# switched back to General right after) to keep the values as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # A
$lastCol  = 51  # AY

$srcRows = @(5, 6, 7)
$dstRows = @(6, 7, 5)   # row 5 -> row 6, row 6 -> row 7, row 7 -> row 5

# Columns that must be forced back to Text after assignment.
$textForceCols = @(9, 25, 27)   # I, Y, AA

# Capture the current contents of rows 5-7 first, because we will
# overwrite them in place.
$capturedValues = @{}

foreach ($r in $srcRows) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $capturedValues["$r-$c"] = $cell.Value()
    }
}

# Write the captured data of each source row into its destination row.
#  - A genuinely blank source cell (.Value() is $null) clears the
#    destination cell, so a populated destination doesn't keep stale data.
#  - A source cell holding an empty string (an explicit-but-empty text
#    cell) is left alone: assigning "" via .Value would turn it into a
#    generic blank cell instead of preserving that empty-text shape, and
#    since those columns are identical (empty) across rows 5-7 to begin
#    with, skipping them is a faithful no-op.
#  - Everything else is written, forcing Text format first for the
#    columns whose strings merely look numeric/date-like.
for ($i = 0; $i -lt $srcRows.Length; $i++) {
    $srcRow = $srcRows[$i]
    $dstRow = $dstRows[$i]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $value = $capturedValues["$srcRow-$c"]
        $destCell = $ws.Cells.Item($dstRow, $c)
        if ($null -eq $value) {
            $destCell.ClearContents()
        } elseif ($value -eq "") {
            continue
        } elseif ($textForceCols -contains $c) {
            $destCell.NumberFormat = "@"
            $destCell.Value = $value
            $destCell.NumberFormat = "General"
        } else {
            $destCell.Value = $value
        }
    }
}
